$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the old bold "Play Adelia: The Fortune Wielder for Free" paragraph
#    that used to sit near the end of the document (its content now lives in
#    a new Meta description paragraph inserted at the top - see step 3).
# ---------------------------------------------------------------------------
$lastCount = $d.Paragraphs.Count
$playAgainPara = $d.Paragraphs($lastCount - 1)
$playAgainPara.Range.Delete()

# ---------------------------------------------------------------------------
# 2) Replace the final italic paragraph's text with the new image-prompt
#    copy, keeping its existing (italic) formatting intact. Scope the Find
#    to just this paragraph's Range so it can't touch the (identical, at the
#    time of writing) text that still lives in the soon-to-be-inserted Meta
#    description paragraph up top.
# ---------------------------------------------------------------------------
$newCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($newCount)
$lastPara.Range.Find.Execute("Read our review of Adelia: The Fortune Wielder, the magical-themed online slot game. Play for free with multiple bonuses and symbol upgrades.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Create a cartoon-style feature image for Adelia: The Fortune Wielder that showcases a happy Maya warrior with glasses. The Maya warrior should be holding a magical wand and surrounded by coins and symbols from the game. The background should be misty and mysterious, with ancient stone structures and a hint of magic in the air. The image should be eye-catching and convey the magical world of Adelia while highlighting the Maya warrior as the main character.", `
    2)

# ---------------------------------------------------------------------------
# 3) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs(2)
$metaXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Adelia: The Fortune Wielder, the magical-themed online slot game. Play for free with multiple bonuses and symbol upgrades.</w:t></w:r></w:p>"
$metaPara.Range.InsertXML($metaXml)
